$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the attendance dates as literal text (e.g. "28/07/2022").
# This edit changes the separator from "/" to "-". Excel's smart entry
# would otherwise reinterpret ambiguous strings like "01-08-2022" as real
# dates, so the whole date column is temporarily switched to Text format,
# written via Value2 (no type coercion), then restored to the default
# "Normal" style so no stray number formatting is left on the cells.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$ws.Range("A3").Value2 = "28-07-2022"
$ws.Range("A4").Value2 = "01-08-2022"
$ws.Range("A5").Value2 = "04-08-2022"
$ws.Range("A6").Value2 = "08-08-2022"
$ws.Range("A7").Value2 = "11-08-2022"
$ws.Range("A8").Value2 = "15-08-2022"
$ws.Range("A9").Value2 = "18-08-2022"
$ws.Range("A10").Value2 = "22-08-2022"
$ws.Range("A11").Value2 = "25-08-2022"
$ws.Range("A12").Value2 = "29-08-2022"
$ws.Range("A13").Value2 = "01-09-2022"
$ws.Range("A14").Value2 = "05-09-2022"
$ws.Range("A15").Value2 = "08-09-2022"
$ws.Range("A16").Value2 = "12-09-2022"
$ws.Range("A17").Value2 = "15-09-2022"
$ws.Range("A18").Value2 = "19-09-2022"
$ws.Range("A19").Value2 = "22-09-2022"
$ws.Range("A20").Value2 = "26-09-2022"
$ws.Range("A21").Value2 = "29-09-2022"

$dateRange.Style = "Normal"

# Row 3 attendance counters: one more valid/invalid record was logged.
$ws.Range("D3").Value2 = 1
$ws.Range("G3").Value2 = 1
